# The "Programa" section of the document contains two paragraphs (a
# Portuguese one and an italicized English one) where items 1-9 of the
# syllabus are all crammed together in a single run with no separation
# between "N." and the end of the previous sentence (e.g. "...condutor.2.
# Celulose..."). The edit inserts a manual line break (<w:br/>) right
# before each item marker "2." through "9." in both paragraphs, turning
# each numbered item into its own visual line within the same run/
# paragraph, without altering any of the wording.

$d = $word.ActiveDocument

# Portuguese paragraph: item markers are written as "N. " (number, dot,
# space) e.g. ".2. Celulose", ".3. Hemiceluloses", etc.
for ($i = 2; $i -le 9; $i++) {
    $d.Content.Find.Execute(
        "\.$i\. ",            # FindText (wildcard): literal dot, digit, literal dot, space
        $true,                 # MatchCase
        $false,                # MatchWholeWord
        $true,                 # MatchWildcards
        $false,                # MatchSoundsLike
        $false,                # MatchAllWordForms
        $true,                 # Forward
        1,                      # Wrap (wdFindContinue)
        $false,                # Format
        ".^l$i. ",             # ReplaceWith: keep the dot, add a line break, keep "N. "
        2                       # Replace (wdReplaceAll)
    ) | Out-Null
}

# English (italic) paragraph: item markers are written as "N." (number,
# dot, no following space) e.g. ".2.Cellulose", ".3.Hemicellulose", etc.
for ($i = 2; $i -le 9; $i++) {
    $d.Content.Find.Execute(
        "\.$i\.",              # FindText (wildcard): literal dot, digit, literal dot
        $true,                 # MatchCase
        $false,                # MatchWholeWord
        $true,                 # MatchWildcards
        $false,                # MatchSoundsLike
        $false,                # MatchAllWordForms
        $true,                 # Forward
        1,                      # Wrap (wdFindContinue)
        $false,                # Format
        ".^l$i.",              # ReplaceWith: keep the dot, add a line break, keep "N."
        2                       # Replace (wdReplaceAll)
    ) | Out-Null
}
